$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "clientes" to "productos"
$ws.Name = "productos"

# Update header row to reflect product data instead of client data
$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "Precio"
$ws.Range("D1").Value = "Cantidad"

# Move the active selection to D1, matching the saved cursor position
$ws.Range("D1").Select()
